# Docx writer: Use different style for block quotes in notes.
#
# Add a new paragraph style "Footnote Block Text" (styleId
# "FootnoteBlockText"), based on / followed-by "Footnote Text", mirroring
# the existing "Block Text" style (which is based on "Body Text") so that
# block quotes inside footnotes can be styled independently from block
# quotes in the main body (e.g. given a different font size).

$d = $word.ActiveDocument

$s = $d.Styles.Add("FootnoteBlockText", 1)  # 1 = wdStyleTypeParagraph

$s.NameLocal = "Footnote Block Text"
$s.BaseStyle = "Footnote Text"
$s.NextParagraphStyle = "Footnote Text"
$s.Priority = 9
$s.UnhideWhenUsed = $true
$s.QuickStyle = $true

$s.ParagraphFormat.SpaceBefore = 5
$s.ParagraphFormat.SpaceAfter = 5
$s.ParagraphFormat.FirstLineIndent = 0
$s.ParagraphFormat.LeftIndent = 24
$s.ParagraphFormat.RightIndent = 24

Write-Output "Added style '$($s.NameLocal)' (FootnoteBlockText)"
